$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.362.29'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.565.18'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.006'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '288.74'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.49%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.06'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3338'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07349'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.75%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.007'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.54'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.822'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.782'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.566.09'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001096'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.70%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '88.59'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06657'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.89%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.095'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.98'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.72'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.78%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.351.90'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.367'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.494'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -10.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.76'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '146.25'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.988'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.56'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.741.74'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.989'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9700'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.824'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.72%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.526'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.05%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.392'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +6.06%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Stellar'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08310'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.53%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02428'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2227'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06259'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.295'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.60%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6112'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.62%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.90'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -6.64%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.006'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.81'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.776'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5715'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.38%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.020'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.03%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '125.35'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.215'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07277'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.28%  '
